# Generate Report for handback
#
# Marks the two source files as handed back (in sync with en-US) in both
# the zh-cn and de-de localization-status sheets, filling in the "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for each, and refreshes the rolled-up Status text on the
# Overview sheet (which shares the same underlying text).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Overview sheet: roll up the new status text for both languages.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: handback info for both tracked files.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

# Row 2 (efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3b0232e625b31d4e5c34cc0b99090e354dbd3f4f/e2e/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
) | Out-Null
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4888992ee70378b0cc1e63e95003a6893ad68a7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.zh-cn.xlf",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.zh-cn.xlf"
) | Out-Null
$zhcn.Range("G2").Value = "2016-01-21 02:53:43"

# Row 3 (ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3b0232e625b31d4e5c34cc0b99090e354dbd3f4f/e2e/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
) | Out-Null
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4888992ee70378b0cc1e63e95003a6893ad68a7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.zh-cn.xlf",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.zh-cn.xlf"
) | Out-Null
$zhcn.Range("G3").Value = "2016-01-21 02:53:43"

# ---------------------------------------------------------------------
# 3. de-de sheet: handback info for both tracked files.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

# Row 2 (efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md)
$dede.Hyperlinks.Add(
    $dede.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3b0232e625b31d4e5c34cc0b99090e354dbd3f4f/e2e/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
) | Out-Null
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01b496f2d021ba175e8e3e6ee9752de4a0d254ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.de-de.xlf",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.de-de.xlf"
) | Out-Null
$dede.Range("G2").Value = "2016-01-21 02:54:04"

# Row 3 (ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md)
$dede.Hyperlinks.Add(
    $dede.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/3b0232e625b31d4e5c34cc0b99090e354dbd3f4f/e2e/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
) | Out-Null
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01b496f2d021ba175e8e3e6ee9752de4a0d254ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.de-de.xlf",
    "",
    "",
    "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.de-de.xlf"
) | Out-Null
$dede.Range("G3").Value = "2016-01-21 02:54:04"
